$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "55+31=86",
    "81+18=99",
    "67-43=24",
    "85+13=98",
    "12+72=84",
    "59-17=42",
    "88-64=24",
    "57+17=74",
    "30+32=62",
    "68-9=59",
    "66-26=40",
    "34+4=38",
    "42+20=62",
    "96-89=7",
    "25-11=14",
    "23+30=53",
    "86-75=11",
    "49+40=89",
    "46+6=52",
    "63-21=42",
    "55-9=46",
    "48-36=12",
    "58+5=63",
    "35+55=90",
    "68-43=25",
    "15+81=96",
    "0+44=44",
    "58-58=0",
    "49-3=46",
    "61+25=86",
    "56+33=89",
    "64-15=49",
    "13+56=69",
    "29-23=6",
    "36+53=89",
    "93-21=72",
    "98-45=53",
    "71-3=68",
    "66-23=43",
    "0+57=57",
    "34+10=44",
    "3+19=22",
    "34+3=37",
    "72+2=74",
    "75+20=95",
    "3+93=96",
    "0+94=94",
    "46-18=28",
    "91-57=34",
    "71+20=91",
    "51-19=32",
    "10+20=30",
    "57+30=87",
    "3+33=36",
    "23-8=15",
    "17+49=66",
    "44-21=23",
    "60-58=2",
    "46+51=97",
    "37-28=9",
    "9+9=18",
    "20+2=22",
    "91-56=35",
    "95-40=55",
    "92-15=77",
    "1+75=76",
    "1+64=65",
    "26-9=17",
    "71-39=32",
    "49-4=45",
    "56+12=68",
    "98-50=48",
    "67+1=68",
    "58+2=60",
    "7+11=18",
    "10-1=9",
    "33+26=59",
    "67-24=43",
    "26+3=29",
    "81-49=32",
    "84+12=96",
    "19+20=39",
    "82+6=88",
    "66+31=97",
    "97-35=62",
    "76+14=90",
    "82-19=63",
    "28+66=94",
    "57+35=92",
    "76-11=65",
    "45-28=17",
    "69+22=91",
    "5+24=29",
    "7+45=52",
    "9+5=14",
    "50+12=62",
    "88-3=85",
    "35+37=72",
    "20+79=99",
    "41+3=44"
)

$cols = 5
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = [math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newValues[$i]
}

Write-Host "Done updating $($newValues.Count) cells"
